{"js": "// Office.js (Word JavaScript API) edit script.\n// Body: async (context) => { ... }\n//\n// The document is a daily division-practice worksheet: a date heading\n// paragraph followed by a table of \"A\u00f7B=C, D\" cells. This edit swaps the\n// date and every table cell's division problem for a new day's values.\n// Because every old value below is unique within the document, we can\n// safely do exact text replacements (search + insertText Replace) without\n// touching formatting (font/size stay exactly as they were).\n\nconst replacements = [\n  [\"2025-05-19 Monday\", \"2025-05-20 Tuesday\"],\n  [\"106\u00f74=26, 2\", \"379\u00f78=47, 3\"],\n  [\"266\u00f76=44, 2\", \"794\u00f77=113, 3\"],\n  [\"940\u00f77=134, 2\", \"215\u00f76=35, 5\"],\n  [\"647\u00f79=71, 8\", \"203\u00f72=101, 1\"],\n  [\"871\u00f73=290, 1\", \"176\u00f76=29, 2\"],\n  [\"830\u00f72=415, 0\", \"169\u00f73=56, 1\"],\n  [\"304\u00f79=33, 7\", \"878\u00f73=292, 2\"],\n  [\"837\u00f74=209, 1\", \"326\u00f74=81, 2\"],\n  [\"460\u00f75=92, 0\", \"318\u00f79=35, 3\"],\n  [\"958\u00f78=119, 6\", \"485\u00f76=80, 5\"],\n  [\"664\u00f72=332, 0\", \"193\u00f75=38, 3\"],\n  [\"146\u00f75=29, 1\", \"586\u00f77=83, 5\"],\n  [\"901\u00f79=100, 1\", \"285\u00f78=35, 5\"],\n  [\"681\u00f76=113, 3\", \"451\u00f72=225, 1\"],\n  [\"545\u00f72=272, 1\", \"663\u00f79=73, 6\"],\n  [\"569\u00f72=284, 1\", \"628\u00f77=89, 5\"],\n  [\"745\u00f75=149, 0\", \"321\u00f75=64, 1\"],\n  [\"809\u00f75=161, 4\", \"953\u00f75=190, 3\"],\n  [\"148\u00f79=16, 4\", \"792\u00f75=158, 2\"],\n  [\"280\u00f79=31, 1\", \"577\u00f74=144, 1\"],\n  [\"414\u00f79=46, 0\", \"139\u00f79=15, 4\"],\n  [\"901\u00f76=150, 1\", \"457\u00f75=91, 2\"],\n  [\"317\u00f73=105, 2\", \"704\u00f76=117, 2\"],\n  [\"573\u00f76=95, 3\", \"758\u00f77=108, 2\"],\n  [\"579\u00f77=82, 5\", \"515\u00f79=57, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The document is a daily division-practice worksheet: a date heading\n# paragraph followed by a table of \"A\u00f7B=C, D\" cells. This edit swaps the\n# date and every table cell's division problem for a new day's values.\n# Because every old value below is unique within the document, exact\n# Find/Replace (wdReplaceAll semantics, one hit each) safely retargets the\n# text runs in place without touching formatting (font/size stay exactly\n# as they were).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-05-19 Monday\", \"2025-05-20 Tuesday\"),\n    @(\"106\u00f74=26, 2\", \"379\u00f78=47, 3\"),\n    @(\"266\u00f76=44, 2\", \"794\u00f77=113, 3\"),\n    @(\"940\u00f77=134, 2\", \"215\u00f76=35, 5\"),\n    @(\"647\u00f79=71, 8\", \"203\u00f72=101, 1\"),\n    @(\"871\u00f73=290, 1\", \"176\u00f76=29, 2\"),\n    @(\"830\u00f72=415, 0\", \"169\u00f73=56, 1\"),\n    @(\"304\u00f79=33, 7\", \"878\u00f73=292, 2\"),\n    @(\"837\u00f74=209, 1\", \"326\u00f74=81, 2\"),\n    @(\"460\u00f75=92, 0\", \"318\u00f79=35, 3\"),\n    @(\"958\u00f78=119, 6\", \"485\u00f76=80, 5\"),\n    @(\"664\u00f72=332, 0\", \"193\u00f75=38, 3\"),\n    @(\"146\u00f75=29, 1\", \"586\u00f77=83, 5\"),\n    @(\"901\u00f79=100, 1\", \"285\u00f78=35, 5\"),\n    @(\"681\u00f76=113, 3\", \"451\u00f72=225, 1\"),\n    @(\"545\u00f72=272, 1\", \"663\u00f79=73, 6\"),\n    @(\"569\u00f72=284, 1\", \"628\u00f77=89, 5\"),\n    @(\"745\u00f75=149, 0\", \"321\u00f75=64, 1\"),\n    @(\"809\u00f75=161, 4\", \"953\u00f75=190, 3\"),\n    @(\"148\u00f79=16, 4\", \"792\u00f75=158, 2\"),\n    @(\"280\u00f79=31, 1\", \"577\u00f74=144, 1\"),\n    @(\"414\u00f79=46, 0\", \"139\u00f79=15, 4\"),\n    @(\"901\u00f76=150, 1\", \"457\u00f75=91, 2\"),\n    @(\"317\u00f73=105, 2\", \"704\u00f76=117, 2\"),\n    @(\"573\u00f76=95, 3\", \"758\u00f77=108, 2\"),\n    @(\"579\u00f77=82, 5\", \"515\u00f79=57, 2\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$null, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]$false, [ref]$false, [ref]$null, 2) | Out-Null\n}\n"}
